# Revert "added requirement openpyxl"
#
# This undoes the "added a requirement openpyxl" commit:
#   1. test_file sheet gains a new "sex" column (inserted as column E,
#      pushing the existing "customer_type" column from E to F) with
#      Male/Female/Other values per row.
#   2. Sheet2 loses its extra demo rows 4-7 (back down to a 3-row sample).
#   3. The active/selected sheet switches from test_file to Sheet2, with
#      Sheet2's selection resting on A3.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("test_file")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- 1. Insert the new "sex" column into test_file -------------------------
$ws1.Columns.Item(5).Insert()
$ws1.Range("E1").Value = "sex"

$sex = @{
    2  = "Male"
    3  = "Other"
    4  = "Female"
    5  = "Male"
    6  = "Female"
    7  = "Male"
    8  = "Female"
    9  = "Male"
    10 = "Female"
    11 = "Male"
    12 = "Male"
    13 = "Male"
    14 = "Female"
    15 = "Female"
    16 = "Female"
    17 = "Female"
    18 = "Female"
    19 = "Male"
    20 = "Female"
    21 = "Male"
    22 = "Male"
    23 = "Female"
    24 = "Female"
    25 = "Female"
    26 = "Male"
    27 = "Male"
    28 = "Male"
    29 = "Male"
    30 = "Female"
}

foreach ($row in $sex.Keys) {
    $ws1.Cells.Item($row, 5).Value = $sex[$row]
}

# --- 2. Drop the extra demo rows from Sheet2 --------------------------------
$ws2.Range("A4:A7").EntireRow.Delete()

# --- 3. Make Sheet2 the active sheet / selection ----------------------------
$ws2.Activate()
$ws2.Range("A3").Select()
